# Apply the change: insert a new row at position 161 (pushing current
# rows 161-178 down to 162-179) on the single worksheet.
#
# The new row 161 copies all values from the (then current) row 161
# except columns D (Fecha) and J (Volumen), which get new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 161; everything currently at row 161 and
# below shifts down by one (row 161 -> row 162, ..., row 178 -> row 179).
$ws.Rows.Item(161).Insert()

# The row that used to be 161 is now row 162; copy its untouched values
# into the freshly inserted row 161, then overwrite D and J with the
# new values from the diff.
$srcRow = 162
$dstRow = 161

$ws.Cells.Item($dstRow, 1).Value = $ws.Cells.Item($srcRow, 1).Value2    # A Mercado ID
$ws.Cells.Item($dstRow, 2).Value = $ws.Cells.Item($srcRow, 2).Value2    # B Mercado
$ws.Cells.Item($dstRow, 3).Value = $ws.Cells.Item($srcRow, 3).Value2    # C Region
$ws.Cells.Item($dstRow, 4).Value = 44449                                 # D Fecha (new)
$ws.Cells.Item($dstRow, 4).NumberFormat = $ws.Cells.Item($srcRow, 4).NumberFormat
$ws.Cells.Item($dstRow, 5).Value = $ws.Cells.Item($srcRow, 5).Value2    # E Codreg
$ws.Cells.Item($dstRow, 6).Value = $ws.Cells.Item($srcRow, 6).Value2    # F Categoria ID
$ws.Cells.Item($dstRow, 7).Value = $ws.Cells.Item($srcRow, 7).Value2    # G Categoria
$ws.Cells.Item($dstRow, 8).Value = $ws.Cells.Item($srcRow, 8).Value2    # H Variedad
$ws.Cells.Item($dstRow, 9).Value = $ws.Cells.Item($srcRow, 9).Value2    # I Calidad
$ws.Cells.Item($dstRow, 10).Value = 65                                   # J Volumen (new)
$ws.Cells.Item($dstRow, 11).Value = $ws.Cells.Item($srcRow, 11).Value2  # K Precio minimo
$ws.Cells.Item($dstRow, 12).Value = $ws.Cells.Item($srcRow, 12).Value2  # L Precio maximo
$ws.Cells.Item($dstRow, 13).Value = $ws.Cells.Item($srcRow, 13).Value2  # M Precio promedio ponderado
$ws.Cells.Item($dstRow, 14).Value = $ws.Cells.Item($srcRow, 14).Value2  # N Unidad de comercializacion
$ws.Cells.Item($dstRow, 15).Value = $ws.Cells.Item($srcRow, 15).Value2  # O Origen
$ws.Cells.Item($dstRow, 16).Value = $ws.Cells.Item($srcRow, 16).Value2  # P Precio $/Kg
$ws.Cells.Item($dstRow, 17).Value = $ws.Cells.Item($srcRow, 17).Value2  # Q Kg o Unidades
$ws.Cells.Item($dstRow, 18).Value = $ws.Cells.Item($srcRow, 18).Value2  # R Clasificacion
